# The "TestData" sheet held stale example rows (TC1/TC2 sample login,
# contact-form and order data) wired up with mailto: hyperlinks. Blank
# out the sample rows (2 and 3) while keeping the header row and the
# existing cell formatting/styles intact, and drop the now-orphaned
# hyperlinks that pointed at the removed sample e-mail addresses.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Remove the hyperlinks attached to the sample data cells.
$ws.Hyperlinks.Delete()

# Clear the sample values in rows 2 and 3 (A:H), leaving row 1 headers
# and all cell styles/number formats untouched.
$ws.Range("A2:H3").ClearContents()

# Keep the UI selection on the now-empty data rows.
$ws.Activate()
$ws.Range("A2:H3").Select()
